$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fill the ACCION column (B) for the new rows 11-17 first (values already exist
# in the shared-string table, so order here has no effect on new string ids).
$ws.Cells.Item(11, 2).Value = "PUT"
$ws.Cells.Item(12, 2).Value = "POST"
$ws.Cells.Item(13, 2).Value = "PUT"
$ws.Cells.Item(14, 2).Value = "GET"
$ws.Cells.Item(15, 2).Value = "GET"
$ws.Cells.Item(16, 2).Value = "GET"
$ws.Cells.Item(17, 2).Value = "GET"

# Rows 11-14 were typed normally: URL (C) then DESCRIPCION (D).
$ws.Cells.Item(11, 3).Value = "/productos/actualizar"
$ws.Cells.Item(11, 4).Value = "Actualiza un producto"

$ws.Cells.Item(12, 3).Value = "/categorias/crear"
$ws.Cells.Item(12, 4).Value = "Crea una categoria"

$ws.Cells.Item(13, 3).Value = "/categorias/actualizar"
$ws.Cells.Item(13, 4).Value = "Actualiza una categoria"

$ws.Cells.Item(14, 3).Value = "/categorias/probarRecuperar"
$ws.Cells.Item(14, 4).Value = "Recupera todas las categorias"

# Rows 15-16: the DESCRIPCION (D) cells were typed before their URL (C) cells.
$ws.Cells.Item(15, 4).Value = "Recupera los pedidos con sus respectivos detalles de un determinado proveedor"
$ws.Cells.Item(16, 4).Value = "Recupera el proveedor, recibiendo el id"

$ws.Cells.Item(15, 3).Value = "/pedidos/buscarPorId{idProveedor}"
$ws.Cells.Item(16, 3).Value = "/proveedores/buscarPorId/{idProveedor}"

# Row 17 typed normally again: URL (C) then DESCRIPCION (D).
$ws.Cells.Item(17, 3).Value = "/productos/buscar/{codigo}"
$ws.Cells.Item(17, 4).Value = "Recupera el producto, recbiendo el código"

# Adjust column D width to match new (wider) content
$ws.Columns.Item(4).ColumnWidth = 73.5

# Update selection to reflect the last-edited cell
$ws.Range("D13").Select()

$wb.Save()
